# iahp_panel_1 settings.xlsx - pairwise testing / boxplots core_addons update
#
# Semantic edits (per commit):
#   - do_normalization (B7):        1 -> 0
#   - do_database_injection (B11):  1 -> 0
#   - do_analysis (B13):            0 -> 1
#   - data_subsets (B17):           "Monos_and_DCs, CD4_T, CD8_T, B"
#                                    -> "Monos_and_DCs, CD4_T, CD8_T, B, TCRgd, NK"
#   - cursor/selection left on B16 (grouping_orders row) when the file was saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B13").Value = 1
$ws.Range("B17").Value = "Monos_and_DCs, CD4_T, CD8_T, B, TCRgd, NK"

$ws.Activate()
$ws.Range("B16").Select()
